$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.4
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 1.8
$ws.Range("L3").Value = 4.33
$ws.Range("Z3").Value = 9.5
$ws.Range("AB3").Value = 23
$ws.Range("AJ3").Value = 7
$ws.Range("AK3").Value = 15

# Row 5
$ws.Range("G5").Value = 2.35
$ws.Range("H5").Value = 2.88
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 1.91
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.44
$ws.Range("S5").Value = 5.5
$ws.Range("T5").Value = 1.14
$ws.Range("AB5").Value = 21
$ws.Range("AR5").Value = 4.2
$ws.Range("AS5").Value = 1.23

# Row 6
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("AA6").Value = 9.5
$ws.Range("AC6").Value = 19
$ws.Range("AJ6").Value = 10
$ws.Range("AK6").Value = 23
$ws.Range("AM6").Value = 51

# Row 9
$ws.Range("G9").Value = 2.55
$ws.Range("I9").Value = 2.4
$ws.Range("J9").Value = 3.4
$ws.Range("AK9").Value = 12
$ws.Range("AM9").Value = 23

# Row 10
$ws.Range("G10").Value = 3.3
$ws.Range("I10").Value = 2
$ws.Range("L10").Value = 2.63
$ws.Range("Z10").Value = 19
$ws.Range("AB10").Value = 41
$ws.Range("AC10").Value = 26

# Row 13
$ws.Range("G13").Value = 2.42
$ws.Range("H13").Value = 2.62
$ws.Range("I13").Value = 3.5
$ws.Range("J13").Value = 3.1
$ws.Range("K13").Value = 1.83
$ws.Range("L13").Value = 4.1
$ws.Range("M13").Value = 1.14
$ws.Range("N13").Value = 5
$ws.Range("Q13").Value = 2.62
$ws.Range("S13").Value = 4.7
$ws.Range("T13").Value = 1.15
$ws.Range("Y13").Value = 5.8
$ws.Range("Z13").Value = 10.5
$ws.Range("AA13").Value = 9.75
$ws.Range("AB13").Value = 27
$ws.Range("AC13").Value = 25
$ws.Range("AD13").Value = 45
$ws.Range("AE13").Value = 5
$ws.Range("AF13").Value = 5.3
$ws.Range("AG13").Value = 16.5
$ws.Range("AJ13").Value = 7.7
$ws.Range("AK13").Value = 17.5
$ws.Range("AL13").Value = 12
$ws.Range("AM13").Value = 55
$ws.Range("AN13").Value = 40
$ws.Range("AO13").Value = 55

# Row 14
$ws.Range("G14").Value = 2.38
$ws.Range("I14").Value = 3.3
$ws.Range("J14").Value = 3.2
$ws.Range("W14").Value = 2.2
$ws.Range("X14").Value = 1.62
$ws.Range("AA14").Value = 10
$ws.Range("AK14").Value = 15
$ws.Range("AL14").Value = 13
$ws.Range("AM14").Value = 41
$ws.Range("AN14").Value = 34

# Row 15
$ws.Range("G15").Value = 1.4
$ws.Range("I15").Value = 6.5
$ws.Range("L15").Value = 7
$ws.Range("W15").Value = 2
$ws.Range("X15").Value = 1.75
$ws.Range("Y15").Value = 7
$ws.Range("AD15").Value = 26
$ws.Range("AI15").Value = 401
$ws.Range("AK15").Value = 34
$ws.Range("AL15").Value = 19

# Row 16
$ws.Range("G16").Value = 2.35
$ws.Range("I16").Value = 2.8
$ws.Range("Q16").Value = 1.95
$ws.Range("R16").Value = 1.9
$ws.Range("S16").Value = 3.25
$ws.Range("T16").Value = 1.33
$ws.Range("AJ16").Value = 9.5
$ws.Range("AL16").Value = 11

# Row 17
$ws.Range("Q17").Value = 1.9
$ws.Range("R17").Value = 1.95

# Row 19
$ws.Range("G19").Value = 2.38
$ws.Range("I19").Value = 2.9
$ws.Range("J19").Value = 3.1
$ws.Range("Y19").Value = 7.5
$ws.Range("Z19").Value = 11
$ws.Range("AA19").Value = 9.5
$ws.Range("AK19").Value = 15

# Row 20
$ws.Range("G20").Value = 1.5
$ws.Range("H20").Value = 4.5
$ws.Range("U20").Value = 1.3
$ws.Range("V20").Value = 3.4
$ws.Range("Z20").Value = 8
$ws.Range("AI20").Value = 201
$ws.Range("AJ20").Value = 17

# Row 22
$ws.Range("G22").Value = 1.7
$ws.Range("H22").Value = 3.75
$ws.Range("J22").Value = 2.38
$ws.Range("K22").Value = 2.1
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 10
$ws.Range("O22").Value = 1.36
$ws.Range("P22").Value = 3
$ws.Range("Q22").Value = 2.1
$ws.Range("R22").Value = 1.7
$ws.Range("S22").Value = 3.75
$ws.Range("T22").Value = 1.25
$ws.Range("W22").Value = 2
$ws.Range("X22").Value = 1.75
$ws.Range("Y22").Value = 6
$ws.Range("AE22").Value = 9
$ws.Range("AF22").Value = 7.5
$ws.Range("AH22").Value = 67
$ws.Range("AI22").Value = 451
$ws.Range("AJ22").Value = 11

# Row 23
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 10
$ws.Range("Q23").Value = 2.08
$ws.Range("R23").Value = 1.73

# Row 24
$ws.Range("O24").Value = 1.5
$ws.Range("P24").Value = 2.5
$ws.Range("U24").Value = 1.62
$ws.Range("V24").Value = 2.2
$ws.Range("AM24").Value = 51
$ws.Range("AP24").Value = 2
$ws.Range("AQ24").Value = 1.85

# Row 26
$ws.Range("H26").Value = 3.1
$ws.Range("J26").Value = 3.75
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 1.07
$ws.Range("N26").Value = 9
$ws.Range("O26").Value = 1.4
$ws.Range("P26").Value = 2.75
$ws.Range("Q26").Value = 2.25
$ws.Range("R26").Value = 1.62
$ws.Range("S26").Value = 4
$ws.Range("T26").Value = 1.22
$ws.Range("W26").Value = 1.95
$ws.Range("X26").Value = 1.8
$ws.Range("AE26").Value = 7.5
$ws.Range("AH26").Value = 51
$ws.Range("AI26").Value = 401
$ws.Range("AJ26").Value = 7
$ws.Range("AK26").Value = 11

# Row 29
$ws.Range("O29").Value = 1.2
$ws.Range("P29").Value = 4.33
$ws.Range("S29").Value = 2.63
$ws.Range("T29").Value = 1.44
